$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("invalidCredentialTest")

# Update row 3 with the new invalid credential test data
$ws.Cells.Item(3, 2).Value = "John124"
$ws.Cells.Item(3, 1).Value = "Kign"

# Remove row 4 entirely (shifts rows up / shrinks used range)
$ws.Rows.Item(4).Delete()

# Update the active selection to the row below the last data row
$ws.Rows.Item(4).Select() | Out-Null
